# This script updates the sheet (h1_data_left) so that the underlying
# "Eig" computation results (14 cerebellar-lobule average columns x 52
# subject rows) reflect a freshly re-run eigen-decomposition, per the
# commit message "Update script to run Eig".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column widths (best match to target authoring tool widths) ----
$ws.Columns.Item(1).ColumnWidth = 13.6
$ws.Columns.Item(2).ColumnWidth = 14.6
$ws.Columns.Item(3).ColumnWidth = 13.76
$ws.Columns.Item(4).ColumnWidth = 13.76
$ws.Columns.Item(5).ColumnWidth = 14.6
$ws.Columns.Item(6).ColumnWidth = 13.76
$ws.Columns.Item(7).ColumnWidth = 13.76
$ws.Columns.Item(8).ColumnWidth = 14.6
$ws.Columns.Item(9).ColumnWidth = 13.6
$ws.Columns.Item(10).ColumnWidth = 13.76
$ws.Columns.Item(11).ColumnWidth = 13.6
$ws.Columns.Item(12).ColumnWidth = 15.6
$ws.Columns.Item(13).ColumnWidth = 15.6
$ws.Columns.Item(14).ColumnWidth = 13.6

# ---- Header row (row 1) - column labels ----
$headerArr = New-Object "object[,]" 1,14
$headerArr[0,0] = "IV_avg"
$headerArr[0,1] = "V_avg"
$headerArr[0,2] = "VI_avg"
$headerArr[0,3] = "Crus_I_avg"
$headerArr[0,4] = "Crus_II_avg"
$headerArr[0,5] = "VIIb_avg"
$headerArr[0,6] = "VIIIa_avg"
$headerArr[0,7] = "VIIIb_avg"
$headerArr[0,8] = "IX_avg"
$headerArr[0,9] = "X_avg"
$headerArr[0,10] = "Vermis_VI_avg"
$headerArr[0,11] = "Vermis_VIIIa_avg"
$headerArr[0,12] = "Vermis_VIIIb_avg"
$headerArr[0,13] = "Vermis_IX_avg"
$ws.Range("A1:N1").Value = $headerArr

# ---- Data rows 2-53: updated Eig averages per subject ----
$dataArr = New-Object "object[,]" 52,14
$dataArr[0,0] = -0.0055525621999999997
$dataArr[0,1] = 0.29095517222
$dataArr[0,2] = 0.23512075750749997
$dataArr[0,3] = 0.22830293296500001
$dataArr[0,4] = 0.0046064838500000024
$dataArr[0,5] = 0.058400845532500002
$dataArr[0,6] = -0.05079566549999999
$dataArr[0,7] = 0.098059993224999989
$dataArr[0,8] = 0.20076685655249998
$dataArr[0,9] = 0.06397760302500001
$dataArr[0,10] = 0.30152674980250005
$dataArr[0,11] = 0.035730512689999987
$dataArr[0,12] = -0.35488959575000001
$dataArr[0,13] = 0.44070753740000002
$dataArr[1,0] = 0.006712059900000034
$dataArr[1,1] = -0.2664735668446
$dataArr[1,2] = -0.50091033336000002
$dataArr[1,3] = -0.16342811632804
$dataArr[1,4] = -0.17235778196399998
$dataArr[1,5] = -0.46440184416000002
$dataArr[1,6] = -0.40838058197999993
$dataArr[1,7] = -0.42675202512000004
$dataArr[1,8] = -0.27290913632000002
$dataArr[1,9] = 0.04510062946
$dataArr[1,10] = -0.65055402271999996
$dataArr[1,11] = -0.1984705067866
$dataArr[1,12] = -0.003480115332000011
$dataArr[1,13] = 0.026836513900000015
$dataArr[2,0] = 0.057299818709999971
$dataArr[2,1] = 0.40242171744800004
$dataArr[2,2] = 0.48175763064780003
$dataArr[2,3] = 0.29403245383999999
$dataArr[2,4] = 0.062670968506000027
$dataArr[2,5] = 0.11124082232000002
$dataArr[2,6] = 0.16684133424460001
$dataArr[2,7] = 0.17316620048
$dataArr[2,8] = -0.093673106633999997
$dataArr[2,9] = 0.80122001865999992
$dataArr[2,10] = 0.51243639978400002
$dataArr[2,11] = -0.30419397741999998
$dataArr[2,12] = -0.39512316246000001
$dataArr[2,13] = -0.41239999745599992
$dataArr[3,0] = 0.1636502921744
$dataArr[3,1] = 0.19744601329
$dataArr[3,2] = 0.36933349498020002
$dataArr[3,3] = 0.011032287973160009
$dataArr[3,4] = 0.11539278562000002
$dataArr[3,5] = 0.24232850349400001
$dataArr[3,6] = 0.29187102542800003
$dataArr[3,7] = 0.15355148034040003
$dataArr[3,8] = 0.19006748640400001
$dataArr[3,9] = 0.30614046572999998
$dataArr[3,10] = 0.32733158996
$dataArr[3,11] = 0.23142978176000001
$dataArr[3,12] = 0.54417414269999997
$dataArr[3,13] = -0.17682135803599999
$dataArr[4,0] = 0.19417307014000001
$dataArr[4,1] = 0.16739789919999998
$dataArr[4,2] = 0.043456365707000018
$dataArr[4,3] = 0.21709814382000001
$dataArr[4,4] = 0.22329697686000002
$dataArr[4,5] = 0.13582836142600002
$dataArr[4,6] = -0.02047633362000003
$dataArr[4,7] = 0.036965527301999999
$dataArr[4,8] = -0.020420584860000004
$dataArr[4,9] = 0.13881046239
$dataArr[4,10] = 0.01538450277600005
$dataArr[4,11] = -0.064897933967999991
$dataArr[4,12] = -0.25292097785399997
$dataArr[4,13] = 0.086686240860000013
$dataArr[5,0] = -0.066708875340000015
$dataArr[5,1] = -0.023894340400000019
$dataArr[5,2] = 0.0039149672400000226
$dataArr[5,3] = -0.043893060460000008
$dataArr[5,4] = -0.12942875327999998
$dataArr[5,5] = 0.035595875599999993
$dataArr[5,6] = 0.10152659094200001
$dataArr[5,7] = 0.35513727547399998
$dataArr[5,8] = 0.086956948898000014
$dataArr[5,9] = 0.32451897641999999
$dataArr[5,10] = -0.08547676055399997
$dataArr[5,11] = -0.12785272421599997
$dataArr[5,12] = -0.24873279271999998
$dataArr[5,13] = -0.18618609751700002
$dataArr[6,0] = -0.78983759099999995
$dataArr[6,1] = -0.75825439342000001
$dataArr[6,2] = -0.76675665717999997
$dataArr[6,3] = -0.47513187771999998
$dataArr[6,4] = -0.53623996621999992
$dataArr[6,5] = -0.79836165613999999
$dataArr[6,6] = -0.7543450390119999
$dataArr[6,7] = -0.47298305827800002
$dataArr[6,8] = -0.41951638457999996
$dataArr[6,9] = -0.65827306576
$dataArr[6,10] = -0.7996763989000002
$dataArr[6,11] = -1.2772924563720001
$dataArr[6,12] = -0.67935170593200001
$dataArr[6,13] = -0.87569763420000013
$dataArr[7,0] = 0.097341796480000006
$dataArr[7,1] = 0.056239976487000007
$dataArr[7,2] = 0.036429902567999981
$dataArr[7,3] = -0.23787216382000001
$dataArr[7,4] = -0.32778947993020002
$dataArr[7,5] = -0.096393437452000019
$dataArr[7,6] = -0.16867439367599996
$dataArr[7,7] = -0.13201610002000003
$dataArr[7,8] = -0.27847579068
$dataArr[7,9] = -0.80340689205999993
$dataArr[7,10] = -0.19817133439999998
$dataArr[7,11] = 0.22157048068000001
$dataArr[7,12] = -0.30716862804599998
$dataArr[7,13] = -0.066026698079999988
$dataArr[8,0] = -0.16104134402
$dataArr[8,1] = -0.110813234486
$dataArr[8,2] = 0.054501290980000006
$dataArr[8,3] = 0.557306479194
$dataArr[8,4] = 0.63925078191999996
$dataArr[8,5] = 0.29690813236819996
$dataArr[8,6] = -0.087068882930000008
$dataArr[8,7] = 0.06738870221200001
$dataArr[8,8] = 0.15053906991999999
$dataArr[8,9] = -0.24992974528000006
$dataArr[8,10] = 0.22772271620157597
$dataArr[8,11] = -0.039885019719999991
$dataArr[8,12] = 0.59435474399799992
$dataArr[8,13] = -0.10991013477600001
$dataArr[9,0] = 0.44031766037999998
$dataArr[9,1] = 0.56947838470800005
$dataArr[9,2] = 0.26257491711999992
$dataArr[9,3] = 0.13575951044000001
$dataArr[9,4] = 0.46548119933999998
$dataArr[9,5] = 0.26982998691400001
$dataArr[9,6] = 0.26457434884400005
$dataArr[9,7] = 0.26141065940000002
$dataArr[9,8] = 0.30682632261999998
$dataArr[9,9] = -0.11986786582000004
$dataArr[9,10] = 0.43478792447400005
$dataArr[9,11] = 0.27292252563999997
$dataArr[9,12] = 0.24835690362400001
$dataArr[9,13] = 0.6202634130400001
$dataArr[10,0] = -0.11971773595200001
$dataArr[10,1] = -0.30641594822000001
$dataArr[10,2] = -0.26321772162399998
$dataArr[10,3] = -0.12454250379780001
$dataArr[10,4] = 0.11555243394
$dataArr[10,5] = -0.1698963381844
$dataArr[10,6] = -0.37586941738000001
$dataArr[10,7] = -0.32475988174199999
$dataArr[10,8] = -0.31098718044200002
$dataArr[10,9] = -0.54658627969999996
$dataArr[10,10] = 0.06605741742
$dataArr[10,11] = -0.025422283880000007
$dataArr[10,12] = -0.23665115122000002
$dataArr[10,13] = 0.21404803948000001
$dataArr[11,0] = 0.19271673462
$dataArr[11,1] = -0.026799019316000004
$dataArr[11,2] = 0.12067170863200002
$dataArr[11,3] = 0.21091892074000002
$dataArr[11,4] = -0.024638657600000002
$dataArr[11,5] = 0.00073925542400000153
$dataArr[11,6] = 0.23950784610199999
$dataArr[11,7] = -0.047125646957999992
$dataArr[11,8] = 0.09009143427999998
$dataArr[11,9] = 0.28339917342000004
$dataArr[11,10] = -0.110724959412
$dataArr[11,11] = 0.50422966044000006
$dataArr[11,12] = 0.30504712514999999
$dataArr[11,13] = 0.44593141543999992
$dataArr[12,0] = 0.050983793136000009
$dataArr[12,1] = 0.086462471759999979
$dataArr[12,2] = 0.065696249591999989
$dataArr[12,3] = -0.093119711792000001
$dataArr[12,4] = -0.195407089078
$dataArr[12,5] = 0.050612295599999989
$dataArr[12,6] = -0.19464958800000001
$dataArr[12,7] = -0.031228342620000026
$dataArr[12,8] = 0.031909025954000005
$dataArr[12,9] = 0.22980019555600001
$dataArr[12,10] = 0.44728591835200004
$dataArr[12,11] = 0.11063538493400001
$dataArr[12,12] = 0.37014382942000001
$dataArr[12,13] = -0.36094546663999993
$dataArr[13,0] = 0.012763213520000006
$dataArr[13,1] = -0.037856889900000003
$dataArr[13,2] = -0.40812266757600002
$dataArr[13,3] = -0.26547512367800002
$dataArr[13,4] = -0.27682526594599999
$dataArr[13,5] = -0.046489470193999981
$dataArr[13,6] = -0.25615233035799995
$dataArr[13,7] = -0.20157759734800001
$dataArr[13,8] = -0.28696590925999999
$dataArr[13,9] = 0.044230556871999983
$dataArr[13,10] = -0.3448673974
$dataArr[13,11] = -0.62361898184800002
$dataArr[13,12] = -1.11642990566
$dataArr[13,13] = -0.14973591454200003
$dataArr[14,0] = 0.049884097252000006
$dataArr[14,1] = -0.14258658960999998
$dataArr[14,2] = -0.091277668074000021
$dataArr[14,3] = -0.38997626023999998
$dataArr[14,4] = -0.19686405086600001
$dataArr[14,5] = -0.33540062044000002
$dataArr[14,6] = -0.24367766515000003
$dataArr[14,7] = 0.023912518020000006
$dataArr[14,8] = -0.020981521447999994
$dataArr[14,9] = -0.098412203430000011
$dataArr[14,10] = 0.091461940173399994
$dataArr[14,11] = -0.58103403388800001
$dataArr[14,12] = -0.41022821136000004
$dataArr[14,13] = -0.10215565262
$dataArr[15,0] = -0.39840686619999999
$dataArr[15,1] = -0.55313472563999999
$dataArr[15,2] = -0.37227434079999994
$dataArr[15,3] = -0.039006404780000004
$dataArr[15,4] = 0.010410451808999998
$dataArr[15,5] = -0.20744706272000002
$dataArr[15,6] = -0.43134750509399999
$dataArr[15,7] = -0.35115131026800006
$dataArr[15,8] = -0.31500464579199999
$dataArr[15,9] = -0.88404562941800013
$dataArr[15,10] = -0.18748416970000001
$dataArr[15,11] = -0.56043554627999992
$dataArr[15,12] = -0.31953788050000009
$dataArr[15,13] = -0.67275482889999993
$dataArr[16,0] = 0.18396212200475001
$dataArr[16,1] = 0.039796053604999992
$dataArr[16,2] = 0.25029212764999997
$dataArr[16,3] = 0.028777114074999992
$dataArr[16,4] = -0.0062678139750000084
$dataArr[16,5] = 0.22733706007499999
$dataArr[16,6] = 0.33936263601224997
$dataArr[16,7] = 0.52113473937249999
$dataArr[16,8] = -0.036152251255750008
$dataArr[16,9] = 0.45184933098925001
$dataArr[16,10] = -0.01902255265
$dataArr[16,11] = 0.20600272222499999
$dataArr[16,12] = 0.15245044762499999
$dataArr[16,13] = 0.10916618907749999
$dataArr[17,0] = -0.026224019334999998
$dataArr[17,1] = -0.2835743168
$dataArr[17,2] = 0.0062989619174999983
$dataArr[17,3] = 0.21934978767500002
$dataArr[17,4] = 0.26310386622699999
$dataArr[17,5] = -0.1553641834275
$dataArr[17,6] = -0.14504638127500002
$dataArr[17,7] = 0.017359474257499993
$dataArr[17,8] = -0.1205955265525
$dataArr[17,9] = 0.26023864387500001
$dataArr[17,10] = -0.34632972407499996
$dataArr[17,11] = -0.11638607713025001
$dataArr[17,12] = 0.21438727499500002
$dataArr[17,13] = -0.99301146544999996
$dataArr[18,0] = -0.15951060623999996
$dataArr[18,1] = -0.16872637566740001
$dataArr[18,2] = -0.11209195788
$dataArr[18,3] = 0.085783716982000005
$dataArr[18,4] = 0.45618384587599997
$dataArr[18,5] = 0.25771638893999993
$dataArr[18,6] = -0.24622251907999998
$dataArr[18,7] = -0.40690361363999994
$dataArr[18,8] = 0.044300823846000002
$dataArr[18,9] = -0.0029966194599999828
$dataArr[18,10] = -0.131265468154
$dataArr[18,11] = -0.15910659693599999
$dataArr[18,12] = -0.14120095723199996
$dataArr[18,13] = 0.18399019502599998
$dataArr[19,0] = 0.13604610534
$dataArr[19,1] = 0.086869779130000005
$dataArr[19,2] = 0.20858752913400003
$dataArr[19,3] = 0.13743745542000002
$dataArr[19,4] = -0.1056272887
$dataArr[19,5] = 0.24452834401999998
$dataArr[19,6] = 0.080356472920000005
$dataArr[19,7] = 0.033488895224000018
$dataArr[19,8] = 0.22593329714600002
$dataArr[19,9] = -0.036498154879999987
$dataArr[19,10] = 0.19704534834000001
$dataArr[19,11] = 0.0091245931199999963
$dataArr[19,12] = 0.15575438780800005
$dataArr[19,13] = -0.083438494439999997
$dataArr[20,0] = 0.051337175860000001
$dataArr[20,1] = -0.068665076583000001
$dataArr[20,2] = -0.075587642839999997
$dataArr[20,3] = 0.21929276573000001
$dataArr[20,4] = 0.28678846564999999
$dataArr[20,5] = -0.12536366203249999
$dataArr[20,6] = -0.19975660252825
$dataArr[20,7] = -0.35959733928900006
$dataArr[20,8] = -0.017409267436600002
$dataArr[20,9] = -0.48981355919750003
$dataArr[20,10] = 0.067173462622499985
$dataArr[20,11] = 0.0033735488750000209
$dataArr[20,12] = 0.13589046165000002
$dataArr[20,13] = 0.035415045425000041
$dataArr[21,0] = -0.39896503902500002
$dataArr[21,1] = -0.78981671419999988
$dataArr[21,2] = -0.577299712375
$dataArr[21,3] = -0.51508202272250003
$dataArr[21,4] = -0.68883008423250003
$dataArr[21,5] = -0.38524105339499992
$dataArr[21,6] = -0.33288944717250002
$dataArr[21,7] = -0.45963625668999997
$dataArr[21,8] = -0.33791248149250003
$dataArr[21,9] = -0.38224265347499997
$dataArr[21,10] = -0.50581664422499995
$dataArr[21,11] = -0.51409383897500005
$dataArr[21,12] = -0.98131437907499985
$dataArr[21,13] = -0.8111446026000001
$dataArr[22,0] = -0.12344101927999999
$dataArr[22,1] = -0.34255980813999998
$dataArr[22,2] = -0.042057240539999971
$dataArr[22,3] = -0.165662974234
$dataArr[22,4] = -0.12396365443800002
$dataArr[22,5] = -0.05063117662000001
$dataArr[22,6] = -0.21733636830399999
$dataArr[22,7] = 0.069084943159999995
$dataArr[22,8] = 0.25310870330259999
$dataArr[22,9] = 0.54333117231999994
$dataArr[22,10] = -0.058782964579999986
$dataArr[22,11] = 0.081499705478000006
$dataArr[22,12] = 0.127787227812
$dataArr[22,13] = 0.49580787550000005
$dataArr[23,0] = 0.068805029739999995
$dataArr[23,1] = 0.15559911307999999
$dataArr[23,2] = 0.26260259043200002
$dataArr[23,3] = 0.17951217152000001
$dataArr[23,4] = 0.24503879292999997
$dataArr[23,5] = 0.23388705855
$dataArr[23,6] = 0.169247126484
$dataArr[23,7] = 0.24306807140000003
$dataArr[23,8] = -0.16828511373400001
$dataArr[23,9] = 0.0032499919753999976
$dataArr[23,10] = 0.26774268724
$dataArr[23,11] = 0.12968320020000004
$dataArr[23,12] = -0.28276733791200004
$dataArr[23,13] = -0.60262392997999992
$dataArr[24,0] = 0.34424082662333327
$dataArr[24,1] = 0.075942803499999975
$dataArr[24,2] = -0.016262875133333337
$dataArr[24,3] = -0.4080991025666667
$dataArr[24,4] = -0.12154137020000001
$dataArr[24,5] = 0.37204609573333336
$dataArr[24,6] = 0.38904063136666661
$dataArr[24,7] = 0.41646727216666662
$dataArr[24,8] = 0.19196533637333332
$dataArr[24,9] = 0.45002121724333333
$dataArr[24,10] = 0.36188497356666671
$dataArr[24,11] = 0.15742090541000001
$dataArr[24,12] = 0.52981178552399999
$dataArr[24,13] = -0.42885114010000008
$dataArr[25,0] = -0.03830385741999999
$dataArr[25,1] = 0.11806609075199999
$dataArr[25,2] = 0.069503853795999992
$dataArr[25,3] = 0.21141246978
$dataArr[25,4] = 0.098260661838000016
$dataArr[25,5] = 0.18937000056600001
$dataArr[25,6] = 0.19080360333800001
$dataArr[25,7] = -0.0089253381799999868
$dataArr[25,8] = 0.10061956674880002
$dataArr[25,9] = -0.32388740240999997
$dataArr[25,10] = -0.021027667309999998
$dataArr[25,11] = 0.23579617714000004
$dataArr[25,12] = 0.18144878820000002
$dataArr[25,13] = -0.33697364403999996
$dataArr[26,0] = -0.082456799659999996
$dataArr[26,1] = -0.22167637076000002
$dataArr[26,2] = -0.023160434217800008
$dataArr[26,3] = 0.36281891101399999
$dataArr[26,4] = 0.31496886504099997
$dataArr[26,5] = -0.21164792509799998
$dataArr[26,6] = -0.45733333802000004
$dataArr[26,7] = 0.048144265638000019
$dataArr[26,8] = -0.024701457698000001
$dataArr[26,9] = -0.20752565656000002
$dataArr[26,10] = -0.1586230307
$dataArr[26,11] = -0.21569056232
$dataArr[26,12] = -0.043336058600000005
$dataArr[26,13] = -0.41218984655999991
$dataArr[27,0] = -0.52413804183200008
$dataArr[27,1] = -0.46929943905999999
$dataArr[27,2] = -0.63986609922000004
$dataArr[27,3] = -0.86429427865999986
$dataArr[27,4] = -1.1171161276400001
$dataArr[27,5] = -1.15044620874
$dataArr[27,6] = -0.74108698183999999
$dataArr[27,7] = -0.74929963554000001
$dataArr[27,8] = -0.45503011008000005
$dataArr[27,9] = -0.24425973179999999
$dataArr[27,10] = -0.77016978456200003
$dataArr[27,11] = -0.73649559115999996
$dataArr[27,12] = -0.32460313700000004
$dataArr[27,13] = -0.48951316725999999
$dataArr[28,0] = -0.15890828265200002
$dataArr[28,1] = -0.1401409039
$dataArr[28,2] = -0.0010828395400000134
$dataArr[28,3] = 0.14072570966600001
$dataArr[28,4] = 0.20738004215600001
$dataArr[28,5] = -0.18482349987800001
$dataArr[28,6] = -0.15627078096000005
$dataArr[28,7] = -0.0080495682479999834
$dataArr[28,8] = 0.11031628862799998
$dataArr[28,9] = -0.13902987788519999
$dataArr[28,10] = -0.020128088470000007
$dataArr[28,11] = -0.88778123751999993
$dataArr[28,12] = -0.57268093804800002
$dataArr[28,13] = -0.19487038945999999
$dataArr[29,0] = 0.19361735155999998
$dataArr[29,1] = 0.14769218606399998
$dataArr[29,2] = 0.18553639941999997
$dataArr[29,3] = 0.14274315460799999
$dataArr[29,4] = 0.13529764380999998
$dataArr[29,5] = -0.011495928771200026
$dataArr[29,6] = 0.21586811683999999
$dataArr[29,7] = 0.11379470878359998
$dataArr[29,8] = 0.086789008024000011
$dataArr[29,9] = 0.027293722359999984
$dataArr[29,10] = 0.34124885090400003
$dataArr[29,11] = -0.29163318398000004
$dataArr[29,12] = -0.73639935876000007
$dataArr[29,13] = -0.1521290388
$dataArr[30,0] = 0.039751457909999996
$dataArr[30,1] = 0.17741998622999999
$dataArr[30,2] = 0.26534508459439998
$dataArr[30,3] = 0.28111276076160002
$dataArr[30,4] = 0.10905401071999998
$dataArr[30,5] = 0.083666756742200005
$dataArr[30,6] = 0.40054555644000001
$dataArr[30,7] = -0.082202134479999994
$dataArr[30,8] = 0.11520343498099998
$dataArr[30,9] = 0.12566246161879999
$dataArr[30,10] = 0.15718671863399999
$dataArr[30,11] = -0.088677030318000005
$dataArr[30,12] = -0.12671048033379997
$dataArr[30,13] = 0.45076407006000008
$dataArr[31,0] = 0.1421138055
$dataArr[31,1] = 0.022168104775999998
$dataArr[31,2] = -0.036010257979999995
$dataArr[31,3] = 0.034343026469999996
$dataArr[31,4] = 0.089460105159940018
$dataArr[31,5] = -0.10168519297600001
$dataArr[31,6] = -0.064179111959399987
$dataArr[31,7] = -0.020056694033599991
$dataArr[31,8] = 0.13291113803999999
$dataArr[31,9] = -0.17686006865999998
$dataArr[31,10] = -0.020044888240000013
$dataArr[31,11] = 0.11563032476000004
$dataArr[31,12] = -0.12285067758000001
$dataArr[31,13] = 0.039735559182000002
$dataArr[32,0] = -0.27659032157500002
$dataArr[32,1] = -0.11083411440000002
$dataArr[32,2] = 0.031434611850000002
$dataArr[32,3] = -0.022942923475250006
$dataArr[32,4] = 0.025919785544999997
$dataArr[32,5] = -0.097007908855500008
$dataArr[32,6] = 0.040410820232499992
$dataArr[32,7] = 0.14362417530250002
$dataArr[32,8] = 0.071210265222499999
$dataArr[32,9] = -0.31709768586749992
$dataArr[32,10] = -0.031785291105000002
$dataArr[32,11] = -0.44022149797499999
$dataArr[32,12] = 0.041267157192500001
$dataArr[32,13] = -0.021149458024999998
$dataArr[33,0] = 0.22030067858000005
$dataArr[33,1] = 0.23885199116600001
$dataArr[33,2] = 0.13626963213600002
$dataArr[33,3] = 0.0061773780320000023
$dataArr[33,4] = -0.096183816642000014
$dataArr[33,5] = 0.022299526778000001
$dataArr[33,6] = -0.091150840627999996
$dataArr[33,7] = -0.126495197754
$dataArr[33,8] = 0.041589306709999996
$dataArr[33,9] = -0.16503839772000001
$dataArr[33,10] = 0.41612881087380005
$dataArr[33,11] = -0.15692491132799999
$dataArr[33,12] = -0.67992293368799994
$dataArr[33,13] = -0.51914120987799994
$dataArr[34,0] = 0.2178724410248
$dataArr[34,1] = 0.25254409235999997
$dataArr[34,2] = 0.17298739634000002
$dataArr[34,3] = 0.14822550142999996
$dataArr[34,4] = 0.16241072403599999
$dataArr[34,5] = -0.012724034740000011
$dataArr[34,6] = 0.13866037009799997
$dataArr[34,7] = 0.03224617867514
$dataArr[34,8] = 0.028413692312000017
$dataArr[34,9] = -0.17617546730000003
$dataArr[34,10] = 0.145981332228
$dataArr[34,11] = -0.038255663075999993
$dataArr[34,12] = -0.23482613144959999
$dataArr[34,13] = -0.083491081240000042
$dataArr[35,0] = 0.055823281953999992
$dataArr[35,1] = 0.28784393976
$dataArr[35,2] = 0.29762660426599996
$dataArr[35,3] = 0.085433904233999991
$dataArr[35,4] = 0.12696970399400001
$dataArr[35,5] = 0.021017830102599999
$dataArr[35,6] = 0.0017440713680000109
$dataArr[35,7] = 0.013609853772000002
$dataArr[35,8] = 0.069148399068000008
$dataArr[35,9] = -0.31053093090200001
$dataArr[35,10] = 0.28902458639999995
$dataArr[35,11] = 0.14303447561399998
$dataArr[35,12] = -0.01050165016000002
$dataArr[35,13] = -0.012989167342999952
$dataArr[36,0] = 0.26547368482
$dataArr[36,1] = 0.33002555882599999
$dataArr[36,2] = 0.39711289796799998
$dataArr[36,3] = 0.1409518495
$dataArr[36,4] = -0.010156996295999999
$dataArr[36,5] = 0.14188298764000001
$dataArr[36,6] = 0.29082018060000003
$dataArr[36,7] = 0.16366837049999999
$dataArr[36,8] = -0.18809263966999998
$dataArr[36,9] = -0.066935986731999964
$dataArr[36,10] = 0.25138656329999998
$dataArr[36,11] = 0.10409388907200001
$dataArr[36,12] = 0.077695759900000022
$dataArr[36,13] = -0.5649416077160001
$dataArr[37,0] = -0.11792149789999999
$dataArr[37,1] = -0.13860648834600001
$dataArr[37,2] = -0.29646539179600001
$dataArr[37,3] = -0.61564347632000005
$dataArr[37,4] = -0.47803412476000001
$dataArr[37,5] = -0.59452899683999993
$dataArr[37,6] = -0.44180484136000003
$dataArr[37,7] = -0.35574196073600001
$dataArr[37,8] = 0.052901485407999993
$dataArr[37,9] = 0.13172222821200003
$dataArr[37,10] = 0.42353508196
$dataArr[37,11] = -0.047994138638000006
$dataArr[37,12] = 0.012919045399999995
$dataArr[37,13] = 0.23359724615999994
$dataArr[38,0] = -0.13432310678800002
$dataArr[38,1] = -0.33073134957399997
$dataArr[38,2] = -0.37720403551800002
$dataArr[38,3] = -0.44349737537999995
$dataArr[38,4] = -0.17402471542
$dataArr[38,5] = -0.29591352540000004
$dataArr[38,6] = -0.53446650820000008
$dataArr[38,7] = -0.23839391467999996
$dataArr[38,8] = -0.30349438072779999
$dataArr[38,9] = -0.34971693498
$dataArr[38,10] = -0.10003399310000001
$dataArr[38,11] = -0.69964353719999994
$dataArr[38,12] = -0.70222983564000008
$dataArr[38,13] = -1.0233699269599998
$dataArr[39,0] = -0.060940422719999998
$dataArr[39,1] = 0.010948954360000008
$dataArr[39,2] = -0.017797432740000008
$dataArr[39,3] = -0.14398827421400001
$dataArr[39,4] = -0.021843875159999991
$dataArr[39,5] = -0.085071401035999986
$dataArr[39,6] = 0.0019885943140000029
$dataArr[39,7] = -0.12633565117199999
$dataArr[39,8] = -0.053008901782399996
$dataArr[39,9] = 0.004022894253999992
$dataArr[39,10] = -0.11466856035999998
$dataArr[39,11] = -0.042099251014000007
$dataArr[39,12] = 0.26938952618599998
$dataArr[39,13] = -0.23226468723400001
$dataArr[40,0] = -0.066363270086000001
$dataArr[40,1] = -0.12533537651999999
$dataArr[40,2] = -0.24275949978200004
$dataArr[40,3] = -0.068492480991999991
$dataArr[40,4] = 0.14786219203999998
$dataArr[40,5] = -0.12722665900000002
$dataArr[40,6] = -0.51912860788460002
$dataArr[40,7] = -0.216928465786
$dataArr[40,8] = 0.33642044014200001
$dataArr[40,9] = -0.48417412546200006
$dataArr[40,10] = -0.16486503499400001
$dataArr[40,11] = -0.49052445321999993
$dataArr[40,12] = -0.20883484488600002
$dataArr[40,13] = 0.013201622960000003
$dataArr[41,0] = -0.23155270479999998
$dataArr[41,1] = -0.44912509881399998
$dataArr[41,2] = -0.46582726691999998
$dataArr[41,3] = -0.30950548227999997
$dataArr[41,4] = -0.18268648628
$dataArr[41,5] = -0.59481571720000004
$dataArr[41,6] = -0.72528131347599989
$dataArr[41,7] = -0.41143099856000004
$dataArr[41,8] = 0.062585443854200007
$dataArr[41,9] = -0.68195201563999996
$dataArr[41,10] = -0.1449635205
$dataArr[41,11] = -0.13017927833999998
$dataArr[41,12] = -0.22045762235999997
$dataArr[41,13] = -0.13393596411200001
$dataArr[42,0] = 0.114451419204
$dataArr[42,1] = 0.093705156117999983
$dataArr[42,2] = -0.092483140790000015
$dataArr[42,3] = -0.42746920566000002
$dataArr[42,4] = -0.4766019011799999
$dataArr[42,5] = -0.23734908754000003
$dataArr[42,6] = -0.17391542339600002
$dataArr[42,7] = -0.086386488051799992
$dataArr[42,8] = -0.25222974064000003
$dataArr[42,9] = 0.021177922567999996
$dataArr[42,10] = -0.0073356742999999506
$dataArr[42,11] = -0.23826197640799998
$dataArr[42,12] = -0.35527757590000009
$dataArr[42,13] = -0.12836198555999995
$dataArr[43,0] = -0.48008888406
$dataArr[43,1] = -0.57388697807599998
$dataArr[43,2] = -0.54484625783999996
$dataArr[43,3] = -0.57143709122000008
$dataArr[43,4] = -0.57460178962000008
$dataArr[43,5] = -0.38788260145800002
$dataArr[43,6] = -0.43652643855399997
$dataArr[43,7] = -0.36900906109999998
$dataArr[43,8] = -0.34911159770600003
$dataArr[43,9] = -0.92601910616000005
$dataArr[43,10] = -0.14503809148000002
$dataArr[43,11] = -0.30844428166000004
$dataArr[43,12] = -0.35431268018839995
$dataArr[43,13] = -1.3213164611999999
$dataArr[44,0] = -0.47868232460000004
$dataArr[44,1] = -0.43560380432249995
$dataArr[44,2] = -0.42223425152499994
$dataArr[44,3] = -0.31438839149999998
$dataArr[44,4] = -0.2956964144575
$dataArr[44,5] = -0.26435875909250001
$dataArr[44,6] = -0.90613836512500001
$dataArr[44,7] = -0.31392575170749998
$dataArr[44,8] = -0.29289079017499997
$dataArr[44,9] = -0.79851855780000003
$dataArr[44,10] = -0.46656751931250001
$dataArr[44,11] = -0.67420687227499998
$dataArr[44,12] = -0.68453887449999995
$dataArr[44,13] = -0.14661040542499998
$dataArr[45,0] = -0.10688718236799999
$dataArr[45,1] = 0.098171948076000021
$dataArr[45,2] = 0.38207788968
$dataArr[45,3] = 0.17463179950599997
$dataArr[45,4] = -0.16233194864000003
$dataArr[45,5] = -0.14596951885600004
$dataArr[45,6] = 0.12773090956000002
$dataArr[45,7] = -0.078736260088399998
$dataArr[45,8] = 0.25332083440000003
$dataArr[45,9] = 0.35950465648000002
$dataArr[45,10] = 0.41659691165999996
$dataArr[45,11] = 0.27371689308000002
$dataArr[45,12] = 0.24229415084
$dataArr[45,13] = 0.35350934163999997
$dataArr[46,0] = 0.084686749160000005
$dataArr[46,1] = -0.0032195897774999899
$dataArr[46,2] = 0.42553151700000003
$dataArr[46,3] = 0.28918847697499994
$dataArr[46,4] = 0.300180705825
$dataArr[46,5] = 0.33619476225
$dataArr[46,6] = 0.39264660557499997
$dataArr[46,7] = 0.24099351227499999
$dataArr[46,8] = 0.16436368227499998
$dataArr[46,9] = 0.12585007207500001
$dataArr[46,10] = 0.60610305253000007
$dataArr[46,11] = -0.011444049433750006
$dataArr[46,12] = 0.14516372687499998
$dataArr[46,13] = -0.50853621854999997
$dataArr[47,0] = 0.23378810196499999
$dataArr[47,1] = 0.63038141783750001
$dataArr[47,2] = 0.6716436443249999
$dataArr[47,3] = 0.37842516137499999
$dataArr[47,4] = 0.038916972139999985
$dataArr[47,5] = 0.24989921540750004
$dataArr[47,6] = 0.3077314356625
$dataArr[47,7] = 0.15569380564999999
$dataArr[47,8] = 0.31096998386000002
$dataArr[47,9] = 0.42613282225224997
$dataArr[47,10] = 0.72940318264999993
$dataArr[47,11] = 0.69048838995000006
$dataArr[47,12] = 0.75699399729999994
$dataArr[47,13] = 0.671610777825
$dataArr[48,0] = 0.22583925387499998
$dataArr[48,1] = -0.045850996464999995
$dataArr[48,2] = -0.032581698575000012
$dataArr[48,3] = 0.020421825275
$dataArr[48,4] = 0.37927454594249999
$dataArr[48,5] = 0.10358890722500003
$dataArr[48,6] = 0.089575759325000032
$dataArr[48,7] = 0.038205670025000013
$dataArr[48,8] = 0.26697534307499998
$dataArr[48,9] = -0.34797055653249997
$dataArr[48,10] = 0.041958458199999993
$dataArr[48,11] = -0.19857843749000004
$dataArr[48,12] = -0.045118542397499996
$dataArr[48,13] = -0.26731702224999998
$dataArr[49,0] = 0.88500557827499993
$dataArr[49,1] = 0.97098749249350003
$dataArr[49,2] = 0.77753865952500001
$dataArr[49,3] = 0.30131401154249998
$dataArr[49,4] = 0.12831525945000002
$dataArr[49,5] = 0.74940368998250007
$dataArr[49,6] = 1.0194872614749999
$dataArr[49,7] = 0.98725081537499992
$dataArr[49,8] = 0.23825909142499996
$dataArr[49,9] = 1.17738800965
$dataArr[49,10] = 0.71044086963749997
$dataArr[49,11] = 1.18414257655
$dataArr[49,12] = 0.94292154485000002
$dataArr[49,13] = 0.430044084125
$dataArr[50,0] = -0.46227527362500009
$dataArr[50,1] = -0.23446563367499998
$dataArr[50,2] = -0.060006015999999995
$dataArr[50,3] = -0.20499590840000004
$dataArr[50,4] = -0.33519966420000002
$dataArr[50,5] = -0.33842290217500004
$dataArr[50,6] = -0.41289689069250002
$dataArr[50,7] = -0.22574283753250002
$dataArr[50,8] = -0.26551841725499997
$dataArr[50,9] = -0.37363148260000001
$dataArr[50,10] = 0.038314317175000001
$dataArr[50,11] = -0.24011343029750004
$dataArr[50,12] = -0.073853443700000021
$dataArr[50,13] = 0.91245794007499992
$dataArr[51,0] = 0.059365105000000008
$dataArr[51,1] = 0.078590702999999998
$dataArr[51,2] = 0.18042258909999997
$dataArr[51,3] = -0.73925974613333334
$dataArr[51,4] = -0.41413917319999999
$dataArr[51,5] = -0.085244305313333338
$dataArr[51,6] = 0.32028219326000001
$dataArr[51,7] = -0.077735135066666683
$dataArr[51,8] = -0.20434717764666668
$dataArr[51,9] = -0.15042423639999999
$dataArr[51,10] = 0.61412935246666656
$dataArr[51,11] = 0.01536585324000001
$dataArr[51,12] = -0.95895243889999993
$dataArr[51,13] = -0.15402198552000004
$ws.Range("A2:N53").Value = $dataArr
